$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 45975
$ws.Range("B2").Value = 14.83
$ws.Range("C2").Value = 17.76
$ws.Range("D2").Value = 13.36
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 6.49
$ws.Range("G2").Value = 11.27
$ws.Range("H2").Value = 23.84
$ws.Range("I2").Value = 50.24
$ws.Range("J2").Value = 56.98
$ws.Range("K2").Value = 48.84
$ws.Range("L2").Value = 29.17
$ws.Range("M2").Value = 17.18
$ws.Range("N2").Value = 11.26
$ws.Range("O2").Value = 10.08
$ws.Range("P2").Value = 15.42
$ws.Range("Q2").Value = 23.19
$ws.Range("R2").Value = 42.67
$ws.Range("S2").Value = 59.27
$ws.Range("T2").Value = 73.25
$ws.Range("U2").Value = 77.75
$ws.Range("V2").Value = 75.13
$ws.Range("W2").Value = 61.34
$ws.Range("X2").Value = 48.86
$ws.Range("Y2").Value = 33.59
$ws.Range("Z2").Value = 34.53
$ws.Range("AB2").Value = 63.24
$ws.Range("AD2").Value = 75.5
$ws.Range("AF2").Value = 68.23999999999999

$wb.Save()
